$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$neo4jDataFile = @'
TC06_CDS_Filter_InstrumentModel-Illumina HiSeq 2500_Neo4jData.xlsx
'@
$webDataFile = @'
TC06_CDS_Filter_InstrumentModel-Illumina HiSeq 2500_WebData.xlsx
'@
$participantsQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p, s, collect(distinct samp.sample_id) as samp
RETURN 
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY `Participant ID`LIMIT 100
'@
$samplesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
ORDER By samp.sample_id LIMIT 100
'@
$filesQuery = @'
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,s,f,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN 
    coalesce(f.file_name, '') as `File Name`,
    coalesce(s.study_name, '') as `Study Name`,
    coalesce(s.phs_accession,'') as `Accession`,
    coalesce(p.participant_id,'') as `Participant ID`,
    coalesce(samp.sample_id, '') as `Sample ID`,
    coalesce(f.file_type, '') as `File Type`
ORDER By f.file_name LIMIT 100
'@
$statQuery = @'
MATCH (f:file)
Match (f)<--(g:genomic_info)
WHERE g.instrument_model in ['Illumina HiSeq 2500']
MATCH (f)-->(samp:sample)-->(p:participant)-->(s:study)
WITH p,f, s, collect(distinct samp.sample_id) as samp
RETURN
count(distinct s) AS Studies,
count(distinct p) AS Participants,
count(distinct samp) AS Samples,
count(distinct f) AS Files
'@

$ws.Range("D2").Value = $neo4jDataFile
$ws.Range("D3").Value = $neo4jDataFile
$ws.Range("D4").Value = $neo4jDataFile
$ws.Range("E2").Value = $webDataFile
$ws.Range("E3").Value = $webDataFile
$ws.Range("E4").Value = $webDataFile

$ws.Range("B2").Value = $participantsQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# Adjust column widths for D and E to reflect new (longer) content (closest achievable value)
$ws.Columns.Item(4).ColumnWidth = 94.02213541666667
$ws.Columns.Item(5).ColumnWidth = 92.30729166666667

# Move the active selection to C2
$ws.Range("C2").Select()

